# "incluido lista rpgpc e atualizado planilha de gastos"
# Update the "Fevereiro" (February) expenses sheet: add two new purchase
# rows (Balança digital / Motorola One) with their values, and move the
# active selection to reflect where the user left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fevereiro")

# New expense rows appended below the existing data / totals.
$ws.Range("A11").Value = "Balança digital"
$ws.Range("C11").Value = 76

$ws.Range("A12").Value = "Motorola One"
$ws.Range("C12").Value = 850

# Leave the selection where the author ended up after entering the data.
[void]$ws.Range("C13").Select()
